# Fulfill 36 scientific names in typo_vlookup
# Adds three new occurrence rows (34-36) to the "Occurrences" sheet,
# mirroring the formatting used by the preceding data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Occurrences")

# ---- Copy direct formatting from the last existing data row (33) down
# ---- into the three new rows so fills/fonts/quote-prefix match.
$ws.Range("F33").Copy()
$ws.Range("F34:F36").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I33:J33").Copy()
$ws.Range("I34:J36").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("P33").Copy()
$ws.Range("P34:P36").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---- Fill column-by-column (matches the order new lookup values were
# ---- typed into the sheet: occurrenceIDs, then kingdom, then names, ...)
$ws.Cells.Item(34, 1).Value = "UNCEN-2000NL-HS001-PM001"
$ws.Cells.Item(35, 1).Value = "UNCEN-2000NL-HS001-PM001"
$ws.Cells.Item(36, 1).Value = "UNCEN-2000NL-HS001-PM001"

$ws.Cells.Item(34, 2).Value = "UNCEN-2000NL-HS001-PM001-VE001"
$ws.Cells.Item(35, 2).Value = "UNCEN-2000NL-HS001-PM001-VE002"
$ws.Cells.Item(36, 2).Value = "UNCEN-2000NL-HS001-PM001-VE003"

$ws.Cells.Item(34, 3).Value = "Human Observation"
$ws.Cells.Item(35, 3).Value = "Human Observation"
$ws.Cells.Item(36, 3).Value = "Human Observation"

$ws.Cells.Item(34, 4).Value = "1999-09-11/1999-09-25"
$ws.Cells.Item(35, 4).Value = "1999-09-11/1999-09-25"
$ws.Cells.Item(36, 4).Value = "1999-09-11/1999-09-25"

$ws.Cells.Item(34, 5).Value = "Animalia"
$ws.Cells.Item(35, 5).Value = "Animalia"
$ws.Cells.Item(36, 5).Value = "Animalia"

$ws.Cells.Item(34, 6).Value = "Phalanger orientalis"
$ws.Cells.Item(35, 6).Value = "Phalanger permixtio"
$ws.Cells.Item(36, 6).Value = "Spilocuscus maculatus"

$ws.Cells.Item(34, 7).Value = "Spesies"
$ws.Cells.Item(35, 7).Value = "Spesies"
$ws.Cells.Item(36, 7).Value = "Spesies"

$ws.Cells.Item(34, 9).Value = "'-2.939800"
$ws.Cells.Item(35, 9).Value = "'-2.939800"
$ws.Cells.Item(36, 9).Value = "'-2.939800"

$ws.Cells.Item(34, 10).Value = "'135.720400"
$ws.Cells.Item(35, 10).Value = "'135.720400"
$ws.Cells.Item(36, 10).Value = "'135.720400"

$ws.Cells.Item(34, 11).Value = "WGS84"
$ws.Cells.Item(35, 11).Value = "WGS84"
$ws.Cells.Item(36, 11).Value = "WGS84"

$ws.Cells.Item(34, 12).Value = "ID"
$ws.Cells.Item(35, 12).Value = "ID"
$ws.Cells.Item(36, 12).Value = "ID"

$ws.Cells.Item(34, 13).Value = "?"
$ws.Cells.Item(35, 13).Value = "?"
$ws.Cells.Item(36, 13).Value = "?"

$ws.Cells.Item(34, 16).Value = "Present"
$ws.Cells.Item(35, 16).Value = "Present"
$ws.Cells.Item(36, 16).Value = "Present"

# ---- Select the newly added occurrenceID range and bring "Occurrences" to front,
# ---- matching the saved workbook's active sheet/selection state.
$ws.Range("P33:P36").Select() | Out-Null
$ws.Activate() | Out-Null
